# Update "Capacità di trasmissione MW" sheet with new transmission-abroad
# capacity values (added transmission abroad for other connections).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacità di trasmissione MW")

# Row 2 (NORD)
$ws.Range("C2").Value = 5500
$ws.Range("D2").Value = 2000
$ws.Range("E2").Value = 2000

# Row 3 (CNOR)
$ws.Range("B3").Value = 4700
$ws.Range("D3").Value = 5200
$ws.Range("H3").Value = 400

# Row 4 (CSUD)
$ws.Range("B4").Value = 2000
$ws.Range("C4").Value = 5900
$ws.Range("E4").Value = 6000
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 1700

# Row 5 (SUD)
$ws.Range("B5").Value = 2000
$ws.Range("D5").Value = 8800

# Row 6 (CALA)
$ws.Range("E6").Value = 5500
$ws.Range("G6").Value = 4100

# Row 7 (SICI)
$ws.Range("D7").Value = 1000
$ws.Range("F7").Value = 4000
$ws.Range("H7").Value = 1000

# Row 8 (SARD)
$ws.Range("C8").Value = 400
$ws.Range("D8").Value = 1900
$ws.Range("G8").Value = 1000

# Row 9 (SUD )
$ws.Range("F9").Value = 4200
